$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 45950
$ws.Cells.Item(2, 2).Value = 7181.2081199582
$ws.Cells.Item(2, 3).Value = 6862.02332252529
$ws.Cells.Item(2, 4).Value = 10340
$ws.Cells.Item(2, 5).Value = 10675.993537
$ws.Cells.Item(2, 6).Value = 0.700364148628826

$ws.Cells.Item(3, 1).Value = 45951
$ws.Cells.Item(3, 2).Value = 7218.14816475568
$ws.Cells.Item(3, 3).Value = 6782.82002216089
$ws.Cells.Item(3, 4).Value = 3620
$ws.Cells.Item(3, 5).Value = 10934.286006
$ws.Cells.Item(3, 6).Value = 286.623244308551

$ws.Cells.Item(4, 1).Value = 45952
$ws.Cells.Item(4, 2).Value = 6883.86755509772
$ws.Cells.Item(4, 3).Value = 6516.66578028577
$ws.Cells.Item(4, 4).Value = 3620
$ws.Cells.Item(4, 5).Value = 10350.430235
$ws.Cells.Item(4, 6).Value = 265.134519174502

$ws.Cells.Item(5, 1).Value = 45953
$ws.Cells.Item(5, 2).Value = 6598.15441220513
$ws.Cells.Item(5, 3).Value = 6232.10648528239
$ws.Cells.Item(5, 4).Value = 3620
$ws.Cells.Item(5, 5).Value = 9847.594343
$ws.Cells.Item(5, 6).Value = 244.231100669886

$ws.Cells.Item(6, 1).Value = 45954
$ws.Cells.Item(6, 2).Value = 6345.86074349242
$ws.Cells.Item(6, 3).Value = 5471.91232215109
$ws.Cells.Item(6, 4).Value = 3620
$ws.Cells.Item(6, 5).Value = 9414.698498
$ws.Cells.Item(6, 6).Value = 205.031253194111

$ws.Cells.Item(7, 1).Value = 45955
$ws.Cells.Item(7, 2).Value = 2044.14689306201
$ws.Cells.Item(7, 3).Value = 3974.84852636007
$ws.Cells.Item(7, 4).Value = 3620
$ws.Cells.Item(7, 5).Value = 5157.791021
$ws.Cells.Item(7, 6).Value = 144.520527262419

$ws.Cells.Item(8, 1).Value = 45956
$ws.Cells.Item(8, 2).Value = 1933.9611505005
$ws.Cells.Item(8, 3).Value = 4010.63328857368
$ws.Cells.Item(8, 4).Value = 3620
$ws.Cells.Item(8, 5).Value = 5037.344132
$ws.Cells.Item(8, 6).Value = 145.584011253049

$ws.Cells.Item(9, 1).Value = 45957
$ws.Cells.Item(9, 2).Value = 6757.64823424925
$ws.Cells.Item(9, 3).Value = 6826.95905158856
$ws.Cells.Item(9, 4).Value = 3620
$ws.Cells.Item(9, 5).Value = 10717.736825
$ws.Cells.Item(9, 6).Value = 298.626985097472

$ws.Cells.Item(10, 1).Value = 45958
$ws.Cells.Item(10, 2).Value = 6757.64823424925
$ws.Cells.Item(10, 3).Value = 7257.38833914105
$ws.Cells.Item(10, 4).Value = 3620
$ws.Cells.Item(10, 5).Value = 10717.736825
$ws.Cells.Item(10, 6).Value = 316.561538745492

$ws.Cells.Item(11, 1).Value = 45959
$ws.Cells.Item(11, 2).Value = 6757.64823424925
$ws.Cells.Item(11, 3).Value = 7295.88321261104
$ws.Cells.Item(11, 4).Value = 3620
$ws.Cells.Item(11, 5).Value = 10717.736825
$ws.Cells.Item(11, 6).Value = 318.165491806742

$ws.Cells.Item(12, 1).Value = 45960
$ws.Cells.Item(12, 2).Value = 6757.64823424925
$ws.Cells.Item(12, 3).Value = 7230.66763906616
$ws.Cells.Item(12, 4).Value = 3620
$ws.Cells.Item(12, 5).Value = 10717.736825
$ws.Cells.Item(12, 6).Value = 315.448176242371

$ws.Cells.Item(13, 1).Value = 45961
$ws.Cells.Item(13, 2).Value = 5678.43345278872
$ws.Cells.Item(13, 3).Value = 5805.08656357129
$ws.Cells.Item(13, 4).Value = 3620
$ws.Cells.Item(13, 5).Value = 9536.254309
$ws.Cells.Item(13, 6).Value = 251.787809157607

$ws.Cells.Item(14, 1).Value = 45962
$ws.Cells.Item(14, 2).Value = 2472.8635019535
$ws.Cells.Item(14, 3).Value = 4048.14168334481
$ws.Cells.Item(14, 4).Value = 6652
$ws.Cells.Item(14, 5).Value = 6332.296934
$ws.Cells.Item(14, 6).Value = 52.3156298079711

$ws.Cells.Item(15, 1).Value = 45963
$ws.Cells.Item(15, 2).Value = 2472.8635019535
$ws.Cells.Item(15, 3).Value = 4107.5611910703
$ws.Cells.Item(15, 4).Value = 6652
$ws.Cells.Item(15, 5).Value = 6332.296934
$ws.Cells.Item(15, 6).Value = 54.7914426298668
